# Apply the "Updated cryptos list" data refresh to sheet1 (coin table).
# Column D (Price) cells are plain numeric-looking text (e.g. "29.219.79",
# "0.3930") that must stay as text, so we force text format before writing
# and restore the default "Normal" style afterwards to avoid leaving behind
# stray number formatting on those cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.219.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.58%  "

$ws.Range("E4").Value = "  -1.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3930"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08428"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.247"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.898.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.348"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("E16").Value = "  -1.77%  "

$ws.Range("E17").Value = "  +2.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001110"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06720"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.026"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.217.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.217"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.114.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.432"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.89"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.063"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1044"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.927"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02477"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06615"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.069"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2193"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.102"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6464"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.235"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.16%  "

$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6035"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.675"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.045"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.167"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.13%  "

Write-Output "Updated 97 cells across 50 rows."
